# Weekly data refresh: a new observation is inserted as row 332 ("Coliflor",
# Vega Central Mapocho de Santiago) and every existing row from 332 downward
# shifts down by one (332->333, ..., 421->422). Excel's native row Insert
# reproduces that shift (and the resulting dimension change to A1:R422)
# automatically, so we only need to insert the row and populate the new
# record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 332, pushing old rows 332:421 down to 333:422.
$ws.Rows("332:332").Insert()

# Populate the new row 332 with the latest weekly reading.
$ws.Cells.Item(332, 1).Value = 9
$ws.Cells.Item(332, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(332, 3).Value = "Metropolitana"
$ws.Cells.Item(332, 4).Value = 44508
$ws.Cells.Item(332, 5).Value = 13
$ws.Cells.Item(332, 6).Value = 100112008
$ws.Cells.Item(332, 7).Value = "Coliflor"
$ws.Cells.Item(332, 8).Value = "Sin especificar"
$ws.Cells.Item(332, 9).Value = "Primera"
$ws.Cells.Item(332, 10).Value = 970
$ws.Cells.Item(332, 11).Value = 600
$ws.Cells.Item(332, 12).Value = 700
$ws.Cells.Item(332, 13).Value = 650
$ws.Cells.Item(332, 14).Value = "$/unidad"
$ws.Cells.Item(332, 15).Value = "Región Metropolitana"
$ws.Cells.Item(332, 16).Value = 650
$ws.Cells.Item(332, 17).Value = 1
$ws.Cells.Item(332, 18).Value = "Hortaliza"
